$d = $word.ActiveDocument

# 1) Expand the sentence about the working "tantárgy felvétel" feature so it
#    also mentions the search function and editing of already uploaded data.
$d.Content.Find.Execute(
    "a tantárgy felvétel. ",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "a tantárgy felvétel, a keresésfunkció, illetve a már feltöltött adatok szerkesztése. ",
    2
) | Out-Null

# 2) Clarify that the "next task" is to extend the subject-enrollment feature
#    (rather than the vague "ezt" / "it").
$d.Content.Find.Execute(
    "feladatunk ezt bővíteni",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "feladatunk a felvételt bővíteni",
    2
) | Out-Null

Write-Output "done"
